# Insert a new weekly record as row 15 in the daily price log.
# All existing rows from 15 downward shift down by one (old row 15 becomes
# row 16, ..., old row 141 becomes row 142), and the newly opened row 15
# is populated with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 15..141 down to 16..142, opening up a blank row 15.
$ws.Rows(15).Insert()

# Fill in the new row 15 with this week's observation.
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = "2021-11-03"
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 100112003
$ws.Range("G15").Value = "Ajo"
$ws.Range("H15").Value = "Chino"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 19000
$ws.Range("M15").Value = 18500
$ws.Range("N15").Value = "$/caja 10 kilos"
$ws.Range("O15").Value = "China"
$ws.Range("P15").Value = 1850
$ws.Range("Q15").Value = 10
$ws.Range("R15").Value = "Hortaliza"
